$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 10000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -10350

$ws.Range("H108").Value = 39900
$ws.Range("J108").Value = 39900
$ws.Range("L108").Value = 39900
$ws.Range("N108").Value = -47580

$ws.Range("H129").Value = 1042.6052
$ws.Range("J129").Value = 1153.6061
$ws.Range("L129").Value = 3460.8183
$ws.Range("N129").Value = -13460.8183

$ws.Range("H138").Value = 2468.6768
$ws.Range("I138").Value = 736.0952
$ws.Range("J138").Value = 2935.141
$ws.Range("K138").Value = 2208.2856
$ws.Range("L138").Value = 8805.423000000001
$ws.Range("M138").Value = 2931.7144
$ws.Range("N138").Value = -19085.423

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 648.05554
$ws.Range("I2").Value = 554.7143
$ws.Range("J2").Value = 974.75
$ws.Range("K2").Value = 554.7143
$ws.Range("L2").Value = 974.75
$ws.Range("M2").Value = -441.7143
$ws.Range("N2").Value = -1200.75

$ws.Range("H61").Value = 1827.1
$ws.Range("I61").Value = 1664.75
$ws.Range("J61").Value = 4100
$ws.Range("K61").Value = 1664.75
$ws.Range("L61").Value = 4100
$ws.Range("M61").Value = -1452.75
$ws.Range("N61").Value = -4524

$ws.Range("H116").Value = 648.05554
$ws.Range("I116").Value = 554.7143
$ws.Range("J116").Value = 974.75
$ws.Range("K116").Value = 554.7143
$ws.Range("L116").Value = 974.75
$ws.Range("M116").Value = 1739.2857
$ws.Range("N116").Value = -5562.75

$ws.Range("H136").Value = 1827.1
$ws.Range("I136").Value = 1664.75
$ws.Range("J136").Value = 4100
$ws.Range("K136").Value = 4994.25
$ws.Range("L136").Value = 12300
$ws.Range("M136").Value = -2444.25
$ws.Range("N136").Value = -17400

$ws.Range("H137").Value = 40706.555
$ws.Range("J137").Value = 40706.555
$ws.Range("L137").Value = 40706.555
$ws.Range("N137").Value = -50906.555

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 648.05554
$ws.Range("I3").Value = 554.7143
$ws.Range("J3").Value = 974.75
$ws.Range("K3").Value = 554.7143
$ws.Range("L3").Value = 974.75
$ws.Range("M3").Value = -440.7143
$ws.Range("N3").Value = -1202.75

$ws.Range("H137").Value = 39522
$ws.Range("J137").Value = 40580
$ws.Range("L137").Value = 40580
$ws.Range("N137").Value = -50780

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2035.6833
$ws.Range("I58").Value = 1741.8302
$ws.Range("J58").Value = 4260.5713
$ws.Range("K58").Value = 1741.8302
$ws.Range("L58").Value = 4260.5713
$ws.Range("M58").Value = -1538.8302
$ws.Range("N58").Value = -4666.5713

$ws.Range("H94").Value = 1556.6923
$ws.Range("I94").Value = 804.8
$ws.Range("J94").Value = 2026.625
$ws.Range("K94").Value = 804.8
$ws.Range("L94").Value = 2026.625
$ws.Range("M94").Value = -353.8
$ws.Range("N94").Value = -2928.625

$ws.Range("H122").Value = 2036.5652
$ws.Range("I122").Value = 1396.25
$ws.Range("J122").Value = 3500.1428
$ws.Range("K122").Value = 4188.75
$ws.Range("L122").Value = 10500.4284
$ws.Range("M122").Value = -1738.75
$ws.Range("N122").Value = -15400.4284

$ws.Range("H132").Value = 4304.1055
$ws.Range("I132").Value = 2160
$ws.Range("J132").Value = 5069.857
$ws.Range("K132").Value = 6480
$ws.Range("L132").Value = 15209.571
$ws.Range("M132").Value = -3950
$ws.Range("N132").Value = -20269.571

$ws.Range("H136").Value = 2035.6833
$ws.Range("I136").Value = 1741.8302
$ws.Range("J136").Value = 4260.5713
$ws.Range("K136").Value = 5225.4906
$ws.Range("L136").Value = 12781.7139
$ws.Range("M136").Value = -2675.4906
$ws.Range("N136").Value = -17881.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1984215.1
$ws.Range("I2").Value = 62
$ws.Range("J2").Value = 2381045.8
$ws.Range("K2").Value = 372
$ws.Range("L2").Value = 14286274.8
$ws.Range("M2").Value = -259
$ws.Range("N2").Value = -14286500.8

$ws.Range("H38").Value = 130.16667
$ws.Range("I38").Value = 99
$ws.Range("J38").Value = 136.4
$ws.Range("K38").Value = 297
$ws.Range("L38").Value = 409.2
$ws.Range("M38").Value = 50
$ws.Range("N38").Value = -1103.2

$ws.Range("H131").Value = 8929603
$ws.Range("J131").Value = 922.5
$ws.Range("L131").Value = 2767.5
$ws.Range("N131").Value = -12847.5

$ws.Range("H137").Value = 2348.5715
$ws.Range("J137").Value = 2781.7646
$ws.Range("L137").Value = 8345.293799999999
$ws.Range("N137").Value = -18545.2938

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 27780254
$ws.Range("I80").Value = 62501850
$ws.Range("J80").Value = 2978
$ws.Range("K80").Value = 62501850
$ws.Range("L80").Value = 2978
$ws.Range("M80").Value = -62500852
$ws.Range("N80").Value = -4974

$ws.Range("H83").Value = 27780254
$ws.Range("I83").Value = 62501850
$ws.Range("J83").Value = 2978
$ws.Range("K83").Value = 312509250
$ws.Range("L83").Value = 14890
$ws.Range("M83").Value = -312504258
$ws.Range("N83").Value = -24874

$ws.Range("H102").Value = 2289.9443
$ws.Range("I102").Value = 1441.6154
$ws.Range("K102").Value = 1441.6154
$ws.Range("M102").Value = 180.3846000000001

$ws.Range("H123").Value = 10889.571
$ws.Range("J123").Value = 10889.571
$ws.Range("L123").Value = 10889.571
$ws.Range("N123").Value = -15789.571

$ws.Range("H137").Value = 61096.43
$ws.Range("J137").Value = 62103.848
$ws.Range("L137").Value = 62103.848
$ws.Range("N137").Value = -72303.848

$ws.Range("H140").Value = 43404
$ws.Range("J140").Value = 43404
$ws.Range("L140").Value = 43404
$ws.Range("N140").Value = -53764

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5999.8335
$ws.Range("I7").Value = 2667.6667
$ws.Range("J7").Value = 7110.5557
$ws.Range("K7").Value = 2667.6667
$ws.Range("L7").Value = 7110.5557
$ws.Range("M7").Value = -2555.6667
$ws.Range("N7").Value = -7334.5557

$ws.Range("H40").Value = 10233.733
$ws.Range("I40").Value = 10521.2
$ws.Range("K40").Value = 10521.2
$ws.Range("M40").Value = -10385.2

$ws.Range("H126").Value = 5999.8335
$ws.Range("I126").Value = 2667.6667
$ws.Range("J126").Value = 7110.5557
$ws.Range("K126").Value = 8003.000100000001
$ws.Range("L126").Value = 21331.6671
$ws.Range("M126").Value = -5533.000100000001
$ws.Range("N126").Value = -26271.6671

$ws.Range("H139").Value = 43288
$ws.Range("J139").Value = 43288
$ws.Range("L139").Value = 43288
$ws.Range("N139").Value = -53568

$ws.Range("H140").Value = 55571.4
$ws.Range("J140").Value = 55571.4
$ws.Range("L140").Value = 55571.4
$ws.Range("N140").Value = -65931.39999999999

$ws.Range("H141").Value = 40204.285
$ws.Range("J141").Value = 40204.285
$ws.Range("L141").Value = 40204.285
$ws.Range("N141").Value = -50564.285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2519.7
$ws.Range("I126").Value = 2105.4666
$ws.Range("K126").Value = 6316.399800000001
$ws.Range("M126").Value = -3846.399800000001

$ws.Range("H132").Value = 7577394
$ws.Range("I132").Value = 970.03845
$ws.Range("J132").Value = 18521118
$ws.Range("K132").Value = 2910.11535
$ws.Range("L132").Value = 55563354
$ws.Range("M132").Value = -380.11535
$ws.Range("N132").Value = -55568414

$ws.Range("H138").Value = 47043.75
$ws.Range("J138").Value = 47043.75
$ws.Range("L138").Value = 47043.75
$ws.Range("N138").Value = -57323.75

$ws.Range("H139").Value = 40013.75
$ws.Range("J139").Value = 39971.332
$ws.Range("L139").Value = 39971.332
$ws.Range("N139").Value = -50251.332

$ws.Range("H140").Value = 27685.8
$ws.Range("J140").Value = 27685.8
$ws.Range("L140").Value = 27685.8
$ws.Range("N140").Value = -38045.8

$ws.Range("H141").Value = 43209.668
$ws.Range("J141").Value = 43209.668
$ws.Range("L141").Value = 43209.668
$ws.Range("N141").Value = -53569.668
